$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change NUMERO_DE_PAGINAS values for rows 2 and 3 from "todas" to numeric 3
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3

# Add new column D: COLETAR, matching formatting of neighbouring column C cells
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "COLETAR"

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "S"

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "S"

$ws.Range("D2").Value = "S"

$excel.CutCopyMode = $false

# Update address text for Pirenopolis row
$ws.Range("A4").Value = "Pirenópolis rua"

# Update the AutoFilter to include the new column
$ws.AutoFilterMode = $false
$ws.Range("A1:D1").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ENDERECOS!_FilterDatabase") {
        $n.RefersTo = "=ENDERECOS!`$A`$1:`$D`$1"
    }
}

$wb.Save()
